$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Id" column before column A, shifting Name/Category/... right.
$ws.Columns.Item(1).EntireColumn.Insert()

$ws.Range("A1").Value = "Id"
$ws.Range("A2").Value = "679c7f07-6191-4a36-a125-a9f7c7e989c1"

# Remove the second data row (the "ooredoo" / SIM Card sale).
$ws.Rows.Item(3).EntireRow.Delete()

# Update the remaining sale's quantity.
$ws.Range("E2").Value = 3

# Update the remaining sale's date/time, keeping them stored as plain text
# (matching the original inline-string representation rather than Excel
# auto-converting them into date/time serial numbers).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2024-09-10"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "17:28:38"
